$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Row 3: height changes from an auto 135 to an explicit custom
#    height of 120.75 (customHeight="1").
# ------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 120.75

# ------------------------------------------------------------------
# 2) Insert a brand new event row ("e037 - Smoke Depletion Phase")
#    right before the current row 38 ("e050 - Evening Debriefing"),
#    pushing every row from 38 on down by one.
# ------------------------------------------------------------------
$ws.Rows.Item(38).Insert()

$ws.Cells.Item(38, 1).Value = "e037"
$ws.Cells.Item(38, 2).Value = "<Bold>e037 Smoke Depletion Phase</Bold> `n<InlineUIContainer><Button Content='r4.71' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    `n<LineBreak/><LineBreak/>`nDeplete smoke in each zone by converting one white full strength Smoke marker to a gray 1/2 strength Smoke marker. Alternatively, remove 1/2 strength Smoke marker. Refer to `n<InlineUIContainer><Button Content='r18.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> for the smoke rules. Click image to continue with `n<InlineUIContainer><Button Content='r4.72' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.`n<LineBreak/><LineBreak/>`n                                              <InlineUIContainer><Image Name='c111Smoke1' Height='100' Width='100'></Image></InlineUIContainer>"

$ws.Rows.Item(38).RowHeight = 120

# ------------------------------------------------------------------
# 3) Update the view state: the active/selected cell moves from B37
#    to B38, and the window is scrolled down a bit further.
# ------------------------------------------------------------------
[void]$ws.Range("B38").Select()
$excel.ActiveWindow.ScrollRow = 37
